# edit.ps1 -- applies the "Methods to Communicate Improvements" rewrite.
#
# NOTE: $d.Content.Find.Execute(..., Replace:=...) runs Word's normal
# AutoCorrect-as-you-type pass over the replacement text, which silently
# turns straight quotes (") into curly/smart quotes. The source document
# uses straight quotes, so every textual change below is done in two
# steps instead:
#   1. Find.Execute() with NO replacement argument, just to locate the
#      range (this does not mutate the document or trigger AutoCorrect).
#   2. Re-wrap that Start/End as a fresh Range and assign to .Text,
#      which substitutes the text verbatim (quotes, dashes, etc. kept
#      exactly as typed) while preserving the run's character formatting.

$d = $word.ActiveDocument
$dash = [char]0x2013

function Get-FoundRange($searchText) {
    # Returns a fresh Range positioned over the first match of $searchText
    # in the whole document, or $null if not found.
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText)
    if (-not $ok) {
        return $null
    }
    return $d.Range($rng.Start, $rng.End)
}

# ---------------------------------------------------------------------
# 1) "Improvement 1 - Method:" -> "Improvement 1 <dash> Communication Method:"
# ---------------------------------------------------------------------
$r = Get-FoundRange("Improvement 1 - ")
$r.Text = "Improvement 1 " + $dash + " Communication "

# ---------------------------------------------------------------------
# 2) Description (Improvement 1): move the opening "(" after the
#    instructional sentence, in front of "For example".
# ---------------------------------------------------------------------
$old2 = ' (Describe the first method you would use. For example: "I would conduct a brief team meeting to discuss the communication improvements and gather feedback from team members. This method allows for open dialogue and ensures everyone understands the rationale behind the changes.")'
$new2 = ' Describe the first method you would use. (For example: "I would conduct a brief team meeting to discuss the communication improvements and gather feedback from team members. This method allows for open dialogue and ensures everyone understands the rationale behind the changes.")'
$r = Get-FoundRange($old2)
$r.Text = $new2

# ---------------------------------------------------------------------
# 3) Reasoning (Improvement 1): same "(" relocation.
# ---------------------------------------------------------------------
$old3 = ' (Explain why this method is appropriate. For example: "Team meetings provide a valuable opportunity for face-to-face interaction and open discussion, which can help build trust and foster a collaborative environment.")'
$new3 = ' Explain why this method is appropriate. (For example: "Team meetings provide a valuable opportunity for face-to-face interaction and open discussion, which can help build trust and foster a collaborative environment.")'
$r = Get-FoundRange($old3)
$r.Text = $new3

# ---------------------------------------------------------------------
# 4) "Improvement 2 - Method:" -> "Improvement 2 <dash> Communication Method:"
# ---------------------------------------------------------------------
$r = Get-FoundRange("Improvement 2 - ")
$r.Text = "Improvement 2 " + $dash + " Communication "

# ---------------------------------------------------------------------
# 5) Description (Improvement 2): reworded to ask for a *second*,
#    different method (drop the reused "first method" example text).
# ---------------------------------------------------------------------
$old5 = ' (Describe the first method you would use. For example: "I would conduct a brief team meeting to discuss the communication improvements and gather feedback from team members. This method allows for open dialogue and ensures everyone understands the rationale behind the changes.") '
$new5 = '  Describe the second method you would use.  It should not be the same as the first.'
$r = Get-FoundRange($old5)
$r.Text = $new5

# ---------------------------------------------------------------------
# 6) Reasoning (Improvement 2): drop the "For example" parenthetical.
# ---------------------------------------------------------------------
$old6 = '(Explain why this method is appropriate. For example: "Team meetings provide a valuable opportunity for face-to-face interaction and open discussion, which can help build trust and foster a collaborative environment.") '
$new6 = 'Explain why this method is appropriate.'
$r = Get-FoundRange($old6)
$r.Text = $new6

# ---------------------------------------------------------------------
# 7) Append a red "(Note: sources are not required)." after the B. item.
# ---------------------------------------------------------------------
$r = Get-FoundRange("Acknowledge sources, using in-text citations and references, for content that is quoted, paraphrased, or summarized.")
$insertPos = $r.End

$noteText = "(Note: sources are not required)."

$spaceRng = $d.Range($insertPos, $insertPos)
$spaceRng.InsertAfter("  ")

$noteStart = $insertPos + 2
$noteRng = $d.Range($noteStart, $noteStart)
$noteRng.InsertAfter($noteText)

$coloredRng = $d.Range($noteStart, $noteStart + $noteText.Length)
$coloredRng.Font.Color = 255
